# Append 12 new NBA game rows (rows 687-698) to Sheet1, growing the
# used range from A1:I686 to A1:I698, and update the selection/scroll
# position to match the final saved view (scrolled to row 665, with
# A687 as the active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column order: Away team, Away Pts, Home team, Home Pts, Overtime,
# Attend., Arena, Win, Loss
$newRows = @(
    @("New York Knicks",       113, "Charlotte Hornets",      92,  "No", 17832, "Spectrum Center",            "New York Knicks",        "Charlotte Hornets"),
    @("Los Angeles Clippers",  108, "Cleveland Cavaliers",    118, "No", 17832, "Rocket Mortgage Fieldhouse",  "Cleveland Cavaliers",    "Los Angeles Clippers"),
    @("New Orleans Pelicans",  112, "Boston Celtics",         118, "No", 17832, "TD Garden",                   "Boston Celtics",         "New Orleans Pelicans"),
    @("Utah Jazz",             114, "Brooklyn Nets",          147, "No", 17832, "Barclays Center",             "Brooklyn Nets",          "Utah Jazz"),
    @("Phoenix Suns",          118, "Miami Heat",             105, "No", 17832, "Kaseya Center",               "Phoenix Suns",           "Miami Heat"),
    @("Los Angeles Lakers",    119, "Houston Rockets",        135, "No", 17832, "Toyota Center",               "Houston Rockets",        "Los Angeles Lakers"),
    @("Sacramento Kings",      103, "Memphis Grizzlies",      94,  "No", 17832, "FedEx Forum",                 "Sacramento Kings",       "Memphis Grizzlies"),
    @("Minnesota Timberwolves",107, "Oklahoma City Thunder",  101, "No", 17832, "Paycom Center",               "Minnesota Timberwolves", "Oklahoma City Thunder"),
    @("Washington Wizards",    118, "San Antonio Spurs",      113, "No", 17832, "Frost Bank Center",           "Washington Wizards",     "San Antonio Spurs"),
    @("Orlando Magic",         129, "Dallas Mavericks",       131, "No", 17832, "American Airlines Center",    "Dallas Mavericks",       "Orlando Magic"),
    @("Milwaukee Bucks",       107, "Denver Nuggets",         113, "No", 17832, "Ball Arena",                  "Denver Nuggets",         "Milwaukee Bucks"),
    @("Philadelphia 76ers",    104, "Portland Trail Blazers", 130, "No", 17832, "Moda Center",                 "Portland Trail Blazers", "Philadelphia 76ers")
)

$startRow = 687
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }

    # Columns B (Away Pts) and D (Home Pts) use a thousands-separator
    # number format in this workbook, same as the preceding rows.
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 4).NumberFormat = "#,##0"
}

# Match the final saved view: scrolled down so row 665 is at the top,
# with A687 selected.
$ws.Range("A687").Select()
$excel.ActiveWindow.ScrollRow = 665
$excel.ActiveWindow.ScrollColumn = 1
